$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# --- New header cells in row 7 (match formatting of existing header cells like A7) ---
$ws.Range("A7").Copy()
$ws.Range("S7:T7").PasteSpecial(-4122)
$ws.Range("S7").Value = "AlarmLoadingDetail"
$ws.Range("T7").Value = "StandbyLoadingDetail"

# --- New data cells in row 8 (match formatting of existing data cell B8) ---
$ws.Range("B8").Copy()
$ws.Range("S8:T8").PasteSpecial(-4122)
$ws.Range("S8").Value = "Battery Alarm (A)"
$ws.Range("T8").Value = "Battery Standby (A)"

# --- New column width for column T ---
$ws.Columns.Item(20).ColumnWidth = 19.6640625

# --- Update sheet view: scroll position and selection ---
$ws.Application.ActiveWindow.ScrollColumn = 13
$ws.Range("P9").Select()
